$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4263.815  # H40: was 4600.1665
$ws.Cells.Item(40, 9).Value = 3865.8235  # I40: was 4100
$ws.Cells.Item(40, 10).Value = 4940.4  # J40: was 5600.5
$ws.Cells.Item(40, 11).Value = 3865.8235  # K40: was 4100
$ws.Cells.Item(40, 12).Value = 4940.4  # L40: was 5600.5
$ws.Cells.Item(40, 13).Value = -3690.8235  # M40: was -3925
$ws.Cells.Item(40, 14).Value = -5290.4  # N40: was -5950.5
$ws.Cells.Item(92, 8).Value = 1546.4  # H92: was 1046.4
$ws.Cells.Item(92, 10).Value = 4783.3335  # J92: was 3116.6667
$ws.Cells.Item(92, 12).Value = 4783.3335  # L92: was 3116.6667
$ws.Cells.Item(92, 14).Value = -7279.3335  # N92: was -5612.6667
$ws.Cells.Item(135, 8).Value = 2152.2222  # H135: was 1345.8667
$ws.Cells.Item(135, 9).Value = 1407.6  # I135: was 904.5833
$ws.Cells.Item(135, 10).Value = 3083  # J135: was 3111
$ws.Cells.Item(135, 11).Value = 12668.4  # K135: was 8141.2497
$ws.Cells.Item(135, 12).Value = 27747  # L135: was 27999
$ws.Cells.Item(135, 13).Value = -10133.4  # M135: was -5606.2497
$ws.Cells.Item(135, 14).Value = -32817  # N135: was -33069
$ws.Cells.Item(138, 8).Value = 2360.4285  # H138: was 2123
$ws.Cells.Item(138, 9).Value = 1284.7  # I138: was 1249.6364
$ws.Cells.Item(138, 10).Value = 5049.75  # J138: was 4524.75
$ws.Cells.Item(138, 11).Value = 3854.1  # K138: was 3748.9092
$ws.Cells.Item(138, 12).Value = 15149.25  # L138: was 13574.25
$ws.Cells.Item(138, 13).Value = 1285.9  # M138: was 1391.0908
$ws.Cells.Item(138, 14).Value = -25429.25  # N138: was -23854.25

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8927.429  # H2: was 7352.933
$ws.Cells.Item(2, 9).Value = 2679.2727  # I2: was 2291.923
$ws.Cells.Item(2, 10).Value = 31837.334  # J2: was 40249.5
$ws.Cells.Item(2, 11).Value = 2679.2727  # K2: was 2291.923
$ws.Cells.Item(2, 12).Value = 31837.334  # L2: was 40249.5
$ws.Cells.Item(2, 13).Value = -2566.2727  # M2: was -2178.923
$ws.Cells.Item(2, 14).Value = -32063.334  # N2: was -40475.5
$ws.Cells.Item(102, 8).Value = 5492.636  # H102: was 5941.9
$ws.Cells.Item(102, 9).Value = 1644  # I102: was 1805
$ws.Cells.Item(102, 11).Value = 1644  # K102: was 1805
$ws.Cells.Item(102, 13).Value = -22  # M102: was -183
$ws.Cells.Item(116, 8).Value = 8927.429  # H116: was 7352.933
$ws.Cells.Item(116, 9).Value = 2679.2727  # I116: was 2291.923
$ws.Cells.Item(116, 10).Value = 31837.334  # J116: was 40249.5
$ws.Cells.Item(116, 11).Value = 2679.2727  # K116: was 2291.923
$ws.Cells.Item(116, 12).Value = 31837.334  # L116: was 40249.5
$ws.Cells.Item(116, 13).Value = -385.2727  # M116: was 2.077000000000226
$ws.Cells.Item(116, 14).Value = -36425.334  # N116: was -44837.5

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8927.429  # H3: was 7352.933
$ws.Cells.Item(3, 9).Value = 2679.2727  # I3: was 2291.923
$ws.Cells.Item(3, 10).Value = 31837.334  # J3: was 40249.5
$ws.Cells.Item(3, 11).Value = 2679.2727  # K3: was 2291.923
$ws.Cells.Item(3, 12).Value = 31837.334  # L3: was 40249.5
$ws.Cells.Item(3, 13).Value = -2565.2727  # M3: was -2177.923
$ws.Cells.Item(3, 14).Value = -32065.334  # N3: was -40477.5
$ws.Cells.Item(105, 8).Value = 1235.75  # H105: was 1303.5454
$ws.Cells.Item(105, 9).Value = 1257.2727  # I105: was 1334
$ws.Cells.Item(105, 11).Value = 1257.2727  # K105: was 1334
$ws.Cells.Item(105, 13).Value = 489.7273  # M105: was 413

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 55.636364  # H2: was 47.583332
$ws.Cells.Item(2, 9).Value = 62.142857  # I2: was 49.333332
$ws.Cells.Item(2, 10).Value = 44.25  # J2: was 42.333332
$ws.Cells.Item(2, 11).Value = 372.857142  # K2: was 295.999992
$ws.Cells.Item(2, 12).Value = 265.5  # L2: was 253.999992
$ws.Cells.Item(2, 13).Value = -259.857142  # M2: was -182.999992
$ws.Cells.Item(2, 14).Value = -491.5  # N2: was -479.999992
$ws.Cells.Item(4, 8).Value = 43565856  # H4: was 71572190
$ws.Cells.Item(4, 10).Value = 794.9231  # J4: was 1599.5
$ws.Cells.Item(4, 12).Value = 2384.7693  # L4: was 4798.5
$ws.Cells.Item(4, 14).Value = -2608.7693  # N4: was -5022.5
$ws.Cells.Item(9, 8).Value = 66.666664  # H9: was 91.25
$ws.Cells.Item(9, 9).Value = 0  # I9: was 92.5
$ws.Cells.Item(9, 10).Value = 66.666664  # J9: was 90
$ws.Cells.Item(9, 11).Value = 0  # K9: was 277.5
$ws.Cells.Item(9, 12).Value = 199.999992  # L9: was 270
$ws.Cells.Item(9, 13).ClearContents()  # M9: was -53.5
$ws.Cells.Item(9, 14).Value = -647.999992  # N9: was -718
$ws.Cells.Item(10, 8).Value = 43.8  # H10: was 52.333332
$ws.Cells.Item(10, 9).Value = 38.333332  # I10: was 52.5
$ws.Cells.Item(10, 11).Value = 114.999996  # K10: was 157.5
$ws.Cells.Item(10, 13).Value = 24.000004  # M10: was -18.5
$ws.Cells.Item(11, 8).Value = 188.5  # H11: was 222.5
$ws.Cells.Item(11, 9).Value = 174.70589  # I11: was 202.22223
$ws.Cells.Item(11, 10).Value = 266.66666  # J11: was 283.33334
$ws.Cells.Item(11, 11).Value = 524.1176700000001  # K11: was 606.66669
$ws.Cells.Item(11, 12).Value = 799.9999799999999  # L11: was 850.0000200000001
$ws.Cells.Item(11, 13).Value = -384.1176700000001  # M11: was -466.66669
$ws.Cells.Item(11, 14).Value = -1079.99998  # N11: was -1130.00002
$ws.Cells.Item(14, 8).Value = 139.3  # H14: was 139.4
$ws.Cells.Item(14, 9).Value = 139.3  # I14: was 139.4
$ws.Cells.Item(14, 11).Value = 417.9  # K14: was 418.2
$ws.Cells.Item(14, 13).Value = -244.9  # M14: was -245.2
$ws.Cells.Item(17, 8).Value = 480.22223  # H17: was 461.94736
$ws.Cells.Item(17, 9).Value = 85.8  # I17: was 84.36364
$ws.Cells.Item(17, 10).Value = 973.25  # J17: was 981.125
$ws.Cells.Item(17, 11).Value = 257.4  # K17: was 253.09092
$ws.Cells.Item(17, 12).Value = 2919.75  # L17: was 2943.375
$ws.Cells.Item(17, 13).Value = -88.39999999999998  # M17: was -84.09092000000001
$ws.Cells.Item(17, 14).Value = -3257.75  # N17: was -3281.375
$ws.Cells.Item(34, 8).Value = 2994  # H34: was 2664.1428
$ws.Cells.Item(34, 10).Value = 3492.8  # J34: was 3024.8333
$ws.Cells.Item(34, 12).Value = 10478.4  # L34: was 9074.499899999999
$ws.Cells.Item(34, 14).Value = -10646.4  # N34: was -9242.499899999999
$ws.Cells.Item(38, 8).Value = 471.43182  # H38: was 452.13043
$ws.Cells.Item(38, 9).Value = 440.07318  # I38: was 430.3095
$ws.Cells.Item(38, 10).Value = 900  # J38: was 681.25
$ws.Cells.Item(38, 11).Value = 1320.21954  # K38: was 1290.9285
$ws.Cells.Item(38, 12).Value = 2700  # L38: was 2043.75
$ws.Cells.Item(38, 13).Value = -973.2195400000001  # M38: was -943.9285
$ws.Cells.Item(38, 14).Value = -3394  # N38: was -2737.75
$ws.Cells.Item(39, 8).Value = 7033.3335  # H39: was 5475
$ws.Cells.Item(39, 9).Value = 0  # I39: was 800
$ws.Cells.Item(39, 11).Value = 0  # K39: was 2400
$ws.Cells.Item(39, 13).ClearContents()  # M39: was -2106
$ws.Cells.Item(40, 8).Value = 51.5  # H40: was 51
$ws.Cells.Item(40, 9).Value = 68.5  # I40: was 51.25
$ws.Cells.Item(40, 10).Value = 34.5  # J40: was 50
$ws.Cells.Item(40, 11).Value = 274  # K40: was 205
$ws.Cells.Item(40, 12).Value = 138  # L40: was 200
$ws.Cells.Item(40, 13).Value = -205  # M40: was -136
$ws.Cells.Item(40, 14).Value = -276  # N40: was -338
$ws.Cells.Item(44, 8).Value = 409.625  # H44: was 471.7143
$ws.Cells.Item(44, 9).Value = 212.83333  # I44: was 263
$ws.Cells.Item(44, 10).Value = 1000  # J44: was 750
$ws.Cells.Item(44, 11).Value = 638.49999  # K44: was 789
$ws.Cells.Item(44, 12).Value = 3000  # L44: was 2250
$ws.Cells.Item(44, 13).Value = -240.49999  # M44: was -391
$ws.Cells.Item(44, 14).Value = -3796  # N44: was -3046
$ws.Cells.Item(46, 8).Value = 523.5  # H46: was 606
$ws.Cells.Item(46, 9).Value = 143  # I46: was 96.5
$ws.Cells.Item(46, 10).Value = 904  # J46: was 809.8
$ws.Cells.Item(46, 11).Value = 429  # K46: was 289.5
$ws.Cells.Item(46, 12).Value = 2712  # L46: was 2429.4
$ws.Cells.Item(46, 13).Value = -338  # M46: was -198.5
$ws.Cells.Item(46, 14).Value = -2894  # N46: was -2611.4
$ws.Cells.Item(51, 8).Value = 4  # H51: was 0
$ws.Cells.Item(51, 9).Value = 4  # I51: was 0
$ws.Cells.Item(51, 11).Value = 12  # K51: was 0
$ws.Cells.Item(51, 13).Value = 448  # M51: was None
$ws.Cells.Item(57, 8).Value = 1731.3334  # H57: was 1875
$ws.Cells.Item(57, 9).Value = 1347  # I57: was 1250
$ws.Cells.Item(57, 11).Value = 4041  # K57: was 3750
$ws.Cells.Item(57, 13).Value = -3482  # M57: was -3191
$ws.Cells.Item(58, 8).Value = 1668  # H58: was 2376.25
$ws.Cells.Item(58, 9).Value = 802.4  # I58: was 1752.5
$ws.Cells.Item(58, 10).Value = 2750  # J58: was 3000
$ws.Cells.Item(58, 11).Value = 2407.2  # K58: was 5257.5
$ws.Cells.Item(58, 12).Value = 8250  # L58: was 9000
$ws.Cells.Item(58, 13).Value = -2279.2  # M58: was -5129.5
$ws.Cells.Item(58, 14).Value = -8506  # N58: was -9256
$ws.Cells.Item(122, 8).Value = 200  # H122: was 921
$ws.Cells.Item(122, 9).Value = 0  # I122: was 832.3333
$ws.Cells.Item(122, 10).Value = 200  # J122: was 1187
$ws.Cells.Item(122, 11).Value = 0  # K122: was 7490.9997
$ws.Cells.Item(122, 12).Value = 1800  # L122: was 10683
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -5040.9997
$ws.Cells.Item(122, 14).Value = -6700  # N122: was -15583

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1676.9286  # H102: was 1139.826
$ws.Cells.Item(102, 9).Value = 1736.6923  # I102: was 1236.75
$ws.Cells.Item(102, 10).Value = 900  # J102: was 493.66666
$ws.Cells.Item(102, 11).Value = 1736.6923  # K102: was 1236.75
$ws.Cells.Item(102, 12).Value = 900  # L102: was 493.66666
$ws.Cells.Item(102, 13).Value = -114.6922999999999  # M102: was 385.25
$ws.Cells.Item(102, 14).Value = -4144  # N102: was -3737.66666
$ws.Cells.Item(107, 8).Value = 394.1  # H107: was 423.77777
$ws.Cells.Item(107, 9).Value = 367.625  # I107: was 419
$ws.Cells.Item(107, 10).Value = 500  # J107: was 433.33334
$ws.Cells.Item(107, 11).Value = 367.625  # K107: was 419
$ws.Cells.Item(107, 12).Value = 500  # L107: was 433.33334
$ws.Cells.Item(107, 13).Value = 1552.375  # M107: was 1501
$ws.Cells.Item(107, 14).Value = -4340  # N107: was -4273.33334
$ws.Cells.Item(126, 8).Value = 6698  # H126: was 6998.1665
$ws.Cells.Item(126, 9).Value = 6316.8  # I126: was 6663
$ws.Cells.Item(126, 11).Value = 18950.4  # K126: was 19989
$ws.Cells.Item(126, 13).Value = -16480.4  # M126: was -17519

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 6250  # H2: was 4266.6665
$ws.Cells.Item(2, 9).Value = 500  # I2: was 400
$ws.Cells.Item(2, 11).Value = 500  # K2: was 400
$ws.Cells.Item(2, 13).Value = -388  # M2: was -288
$ws.Cells.Item(16, 8).Value = 699.5  # H16: was 719
$ws.Cells.Item(16, 9).Value = 699.5  # I16: was 719
$ws.Cells.Item(16, 11).Value = 699.5  # K16: was 719
$ws.Cells.Item(16, 13).Value = -529.5  # M16: was -549
$ws.Cells.Item(82, 8).Value = 6716.6665  # H82: was 3848.4167
$ws.Cells.Item(82, 9).Value = 4150  # I82: was 1772.625
$ws.Cells.Item(82, 11).Value = 4150  # K82: was 1772.625
$ws.Cells.Item(82, 13).Value = -3789  # M82: was -1411.625
$ws.Cells.Item(85, 8).Value = 6716.6665  # H85: was 3848.4167
$ws.Cells.Item(85, 9).Value = 4150  # I85: was 1772.625
$ws.Cells.Item(85, 11).Value = 4150  # K85: was 1772.625
$ws.Cells.Item(85, 13).Value = -2902  # M85: was -524.625

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 843  # H81: was 1034.75
$ws.Cells.Item(81, 9).Value = 843  # I81: was 1034.75
$ws.Cells.Item(81, 11).Value = 1686  # K81: was 2069.5
$ws.Cells.Item(81, 13).Value = -625  # M81: was -1008.5
$ws.Cells.Item(84, 8).Value = 843  # H84: was 1034.75
$ws.Cells.Item(84, 9).Value = 843  # I84: was 1034.75
$ws.Cells.Item(84, 11).Value = 8430  # K84: was 10347.5
$ws.Cells.Item(84, 13).Value = -3126  # M84: was -5043.5
$ws.Cells.Item(96, 8).Value = 2678.4  # H96: was 1856.5714
$ws.Cells.Item(96, 9).Value = 2448  # I96: was 1566
$ws.Cells.Item(96, 11).Value = 2448  # K96: was 1566
$ws.Cells.Item(96, 13).Value = -1075  # M96: was -193

Write-Host "Applied 198 cell updates across 7 sheets"